$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): extend with two new labeled columns P1=14, Q1=15 ---
# Copy formatting from the existing last header cell (O1, col 15) into the
# new cells (P1 col 16, Q1 col 17), then set their values.
$ws.Cells.Item(1, 15).Copy($ws.Cells.Item(1, 16))
$ws.Cells.Item(1, 16).Value = 14
$ws.Cells.Item(1, 15).Copy($ws.Cells.Item(1, 17))
$ws.Cells.Item(1, 17).Value = 15

# --- Data rows 2-25: swap values in columns I/K and M/O, and append P/Q=2 ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P: new column, value 2
    $ws.Cells.Item($r, 17).Value = 2   # Q: new column, value 2
}
